# Master Allocation sheet: update cost_code (column G) values for rows that
# currently read "9000 100F" so that billing rows flag the need for a cost
# code / cross-charge, and clear out two stray leftover summary rows.
#
# Most rows simply get " / CC NEEDED" appended. Three rows (121-123) are a
# special case that also needs the "9000 100M" code folded in. Rows 488-489
# are stray duplicate/summary rows where the cost_code (and for row 488 the
# units/rate/amount too) should be wiped out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple case: "9000 100F" -> "9000 100F / CC NEEDED" -------------------
# (contiguous row blocks across column G where the cost_code was "9000 100F")
$ws.Range("G92:G98").Value   = "9000 100F / CC NEEDED"
$ws.Range("G105:G109").Value = "9000 100F / CC NEEDED"
$ws.Range("G114").Value      = "9000 100F / CC NEEDED"
$ws.Range("G117:G120").Value = "9000 100F / CC NEEDED"
$ws.Range("G124:G127").Value = "9000 100F / CC NEEDED"
$ws.Range("G143:G315").Value = "9000 100F / CC NEEDED"
$ws.Range("G317:G477").Value = "9000 100F / CC NEEDED"
$ws.Range("G486").Value      = "9000 100F / CC NEEDED"

# --- Special case: rows 121-123 also need the 9000 100M code noted ---------
$ws.Range("G121:G123").Value = "9000 100M | 9000 100F | CC NEEDED"

# --- Stray rows 488/489: wipe out leftover placeholder data -----------------
$ws.Range("G488").ClearContents()
$ws.Range("I488").ClearContents()
$ws.Range("K488").ClearContents()
$ws.Range("L488").ClearContents()

$ws.Range("G489").ClearContents()
